$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp in title cell A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 18:52"

# Row 4
$ws.Range("B4").Value = 566654
$ws.Range("C4").Value = 6354
$ws.Range("D4").Value = 33743
$ws.Range("E4").Value = 510034
$ws.Range("G4").Value = 772
$ws.Range("H4").Value = 22877

# Row 8
$ws.Range("B8").Value = 128002
$ws.Range("C8").Value = 148
$ws.Range("E8").Value = 60672
$ws.Range("G8").Value = 8
$ws.Range("H8").Value = 3030

# Row 12
$ws.Range("B12").Value = 61049
$ws.Range("C12").Value = 4093
$ws.Range("D12").Value = 3957
$ws.Range("E12").Value = 55796
$ws.Range("F12").Value = 1786
$ws.Range("G12").Value = 98
$ws.Range("H12").Value = 1296

# Row 14
$ws.Range("F14").Value = 1358

# Row 26
$ws.Range("A26").Value = "Ecuador"
$ws.Range("B26").Value = 7529
$ws.Range("C26").Value = 63
$ws.Range("D26").Value = 597
$ws.Range("E26").Value = 6577
$ws.Range("F26").Value = 121
$ws.Range("G26").Value = 22
$ws.Range("H26").Value = 355

# Row 27
$ws.Range("A27").Value = "Chile"
$ws.Range("B27").Value = 7525
$ws.Range("C27").Value = 312
$ws.Range("D27").Value = 2367
$ws.Range("E27").Value = 5076
$ws.Range("F27").Value = 387
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 82

# Row 28
$ws.Range("A28").Value = "Peru"
$ws.Range("B28").Value = 7519
$ws.Range("D28").Value = 1798
$ws.Range("E28").Value = 5528
$ws.Range("F28").Value = 134
$ws.Range("H28").Value = 193

# Row 36
$ws.Range("B36").Value = 5496
$ws.Range("C36").Value = 266
$ws.Range("E36").Value = 4308

# Row 67
$ws.Range("E67").Value = 1202
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 28

# Row 108
$ws.Range("A108").Value = "Guinea"
$ws.Range("B108").Value = 319
$ws.Range("C108").Value = 69
$ws.Range("D108").Value = 17
$ws.Range("E108").Value = 302
$ws.Range("H108").Value = 0

# Row 109
$ws.Range("A109").Value = "Estado de Palestina"
$ws.Range("B109").Value = 308
$ws.Range("C109").Value = 18
$ws.Range("D109").Value = 58
$ws.Range("E109").Value = 248

# Row 110
$ws.Range("A110").Value = "Republica de Yibuti"
$ws.Range("B110").Value = 298
$ws.Range("C110").Value = 84
$ws.Range("D110").Value = 41
$ws.Range("E110").Value = 255
$ws.Range("F110").Value = 0

# Row 111
$ws.Range("A111").Value = "Senegal"
$ws.Range("B111").Value = 291
$ws.Range("C111").Value = 11
$ws.Range("D111").Value = 178
$ws.Range("E111").Value = 111
$ws.Range("F111").Value = 1
$ws.Range("H111").Value = 2

# Row 112
$ws.Range("A112").Value = "Montenegro"
$ws.Range("B112").Value = 273
$ws.Range("C112").Value = 1
$ws.Range("D112").Value = 5
$ws.Range("E112").Value = 265
$ws.Range("F112").Value = 7

# Row 113
$ws.Range("A113").Value = "Georgia"
$ws.Range("B113").Value = 272
$ws.Range("C113").Value = 15
$ws.Range("D113").Value = 67
$ws.Range("E113").Value = 202
$ws.Range("F113").Value = 6
$ws.Range("H113").Value = 3

# Row 114
$ws.Range("A114").Value = "Vietnam"
$ws.Range("B114").Value = 265
$ws.Range("C114").Value = 3
$ws.Range("D114").Value = 146
$ws.Range("E114").Value = 119
$ws.Range("F114").Value = 8

# Row 154
$ws.Range("B154").Value = 45
$ws.Range("C154").Value = 2
$ws.Range("E154").Value = 13

# Row 163
$ws.Range("E163").Value = 21
$ws.Range("G163").Value = 2
$ws.Range("H163").Value = 4

# Row 167
$ws.Range("A167").Value = "Republica del Chad"
$ws.Range("C167").Value = 5
$ws.Range("D167").Value = 2
$ws.Range("F167").Value = 0
$ws.Range("H167").Value = 0

# Row 168
$ws.Range("A168").Value = "Antigua y Barbuda"
$ws.Range("C168").Value = 2
$ws.Range("D168").Value = 0
$ws.Range("F168").Value = 1
$ws.Range("H168").Value = 2

# Row 182
$ws.Range("A182").Value = "Suazilandia"
$ws.Range("B182").Value = 15
$ws.Range("C182").Value = 1
$ws.Range("D182").Value = 7
$ws.Range("E182").Value = 8
$ws.Range("F182").Value = 0

# Row 183
$ws.Range("A183").Value = "Granada"
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 14
$ws.Range("F183").Value = 2

# Row 184
$ws.Range("A184").Value = "Nepal"
$ws.Range("C184").Value = 2
$ws.Range("D184").Value = 1
$ws.Range("E184").Value = 13
$ws.Range("F184").Value = 0
$ws.Range("H184").Value = 0

# Row 185
$ws.Range("A185").Value = "Belice"
$ws.Range("E185").Value = 12
$ws.Range("F185").Value = 1
$ws.Range("H185").Value = 2

# Row 186
$ws.Range("A186").Value = "Zimbabue"
$ws.Range("D186").Value = 0
$ws.Range("E186").Value = 11
$ws.Range("H186").Value = 3

# Row 205
$ws.Range("A205").Value = "Islas Malvinas"
$ws.Range("D205").Value = 1
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# Row 206
$ws.Range("A206").Value = "Burundi"
$ws.Range("D206").Value = 0
$ws.Range("G206").Value = 1
$ws.Range("H206").Value = 1

# Row 215
$ws.Range("A215").Value = "San Pedro y Miquelon"

# Row 216
$ws.Range("A216").Value = "Yemen"
